$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add a new "CNN" results block, mirroring the layout used by the other
# model blocks above it (Decision Tree / Random Forest / Baysian
# Classifier / SVM), starting at row 76.
# ---------------------------------------------------------------------

# Section title
$ws.Range("A76").Value = "CNN"

# Column headers (Cleveland / Virginia / Hungarian / Switzerland / Centralized / AVG)
$ws.Range("C77").Value = "Cleveland"
$ws.Range("D77").Value = "Virginia"
$ws.Range("E77").Value = "Hungarian"
$ws.Range("F77").Value = "Switzerland"
$ws.Range("G77").Value = "Centralized"
$ws.Range("I77").Value = "AVG"

# Metric rows (values intentionally left blank - not measured yet)
$ws.Range("B78").Value = "Precision - 1"
$ws.Range("B79").Value = "Recall - 1"
$ws.Range("B81").Value = "Precision - 0"
$ws.Range("B82").Value = "Recall - 0"
$ws.Range("B84").Value = "ACC"

# Reuse the "output" cell style (thin border / grey fill) that highlights
# the Centralized column on the other blocks, copying it from the SVM
# block immediately above so the existing style index is reused instead
# of a new one being created.
$ws.Range("G65").Copy()
$ws.Range("G78").PasteSpecial(-4122)
$ws.Range("G79").PasteSpecial(-4122)
$ws.Range("G84").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Reflect the new block in the sheet's selection/scroll position, as a
# user would leave it after entering the new data.
$ws.Range("C78:I84").Select()
